# Correct over-long "valor pago" figures that were accidentally computed
# from the wrong column (customer-name length bug) -- replace each
# mis-computed total with the correctly scaled value.

$d = $word.ActiveDocument

$replacements = @(
    @("27.539,62", "275,40"),
    @("15.496,27", "154,96"),
    @("17.554,22", "175,54"),
    @("9.611,11",  "96,11"),
    @("17.019,02", "170,19"),
    @("32.146,05", "321,46"),
    @("12.444,52", "124,45"),
    @("28.188,30", "281,88"),
    @("19.481,41", "194,81"),
    @("14.910,08", "149,10"),
    @("15.309,34", "153,09"),
    @("16.031,43", "160,31"),
    @("26.175,60", "261,76"),
    @("43.124,40", "431,24"),
    @("16.544,00", "165,44"),
    @("17.309,60", "173,10"),
    @("8.683,41",  "86,83"),
    @("6.783,14",  "67,83"),
    @("55.474,50", "554,74"),
    @("13.536,64", "135,37")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
